# actualizacion graficos y tablas
# Add the 2020 data point to the "Datos" sheet (Fecha=2020, Valor=6.2),
# following the existing Fecha/Valor table that runs from row 1 (headers)
# through row 21 (2019, 6.8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Force the "2020" label to be written as text (matching the other year
# labels in column A), not auto-converted to a number, then restore the
# cell's style to Normal so no stray number-format styling is left behind.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "2020"
$ws.Range("A22").Style = "Normal"

$ws.Range("B22").Value = 6.2
